$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing AutoFilter criteria (State=... / Physical Therapists / 2028 / South)
$ws.AutoFilterMode = $false

# Re-apply AutoFilter over the full data range (now including the totals row 505)
# with the new criteria: State = Arizona AND OccupationName = Physical Therapists
$rng = $ws.Range("A1:O505")
$rng.AutoFilter(2, @("Arizona"))
$rng.AutoFilter(4, @("Physical Therapists"))

# Keep the _FilterDatabase defined name in sync with the new filter range
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$O`$505"

# Match the author's final cell selection
$ws.Range("A205").Select()
